# Prefix each step-row's Name (column A) with the worksheet's own name,
# for every "protocol" worksheet in the workbook (all sheets except the
# first six overview/journey sheets). The header cell "Name" in row 1 is
# left untouched.

$wb = $excel.ActiveWorkbook

# Sheets that must NOT be touched (overview / journey sheets).
$excluded = @("DamonJourney", "MeetupRedirect", "NRWaves", "PersonalDamon", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    if ($excluded -contains $ws.Name) {
        continue
    }

    $prefix = $ws.Name
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($current -ne $null -and $current -ne "") {
            $cell.Value2 = "$prefix $current"
        }
    }
}
